$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 135; existing rows 135-144 shift down to 136-145.
$ws.Rows.Item(135).Insert()

# Populate the new row 135 with the new weekly data point.
$ws.Cells.Item(135, 1).Value = 7
$ws.Cells.Item(135, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(135, 3).Value = "Ñuble"
$ws.Cells.Item(135, 4).Value = 45075
$ws.Cells.Item(135, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(135, 5).Value = 16
$ws.Cells.Item(135, 6).Value = 100112037
$ws.Cells.Item(135, 7).Value = "Cebollín"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 100
$ws.Cells.Item(135, 11).Value = 6000
$ws.Cells.Item(135, 12).Value = 6000
$ws.Cells.Item(135, 13).Value = 6000
$ws.Cells.Item(135, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(135, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(135, 16).Value = 167
$ws.Cells.Item(135, 17).Value = 36
$ws.Cells.Item(135, 18).Value = "Hortaliza"
